$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.739.93'
$ws.Range("E2").Value = '  +1.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.453.97'
$ws.Range("E3").Value = '  +1.46%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.18'
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.26'
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +15.49%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.461.71'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("E10").Value = '  -1.46%  '
$ws.Range("E11").Value = '  +2.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  +3.71%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.059.65'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.17'
$ws.Range("E16").Value = '  +3.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.782.46'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.496.55'
$ws.Range("E18").Value = '  +2.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("E19").Value = '  +3.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.33'
$ws.Range("E20").Value = '  +2.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '380.38'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.09'
$ws.Range("E22").Value = '  +0.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.551'
$ws.Range("E23").Value = '  +3.79%  '
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.79'
$ws.Range("E24").Value = '  +1.43%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.00'
$ws.Range("E27").Value = '  +6.65%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.52'
$ws.Range("E30").Value = '  +12.03%  '
$ws.Range("E31").Value = '  +0.07%  '
$ws.Range("E32").Value = '  +2.15%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.55'
$ws.Range("E33").Value = '  +1.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.21'
$ws.Range("E34").Value = '  +6.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("E35").Value = '  +11.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.47'
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  +5.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0774'
$ws.Range("E38").Value = '  +2.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.951.06'
$ws.Range("E39").Value = '  -0.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '26.66'
$ws.Range("E40").Value = '  -0.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.59'
$ws.Range("E41").Value = '  +7.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.55'
$ws.Range("E42").Value = '  +2.40%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.93'
$ws.Range("E43").Value = '  +2.74%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0318'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.18'
$ws.Range("E45").Value = '  +13.13%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.775'
$ws.Range("E46").Value = '  +2.39%  '
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.112'
$ws.Range("E48").Value = '  +8.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '310.09'
$ws.Range("E49").Value = '  +5.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.872'
$ws.Range("E50").Value = '  +4.57%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.18'
$ws.Range("E51").Value = '  -0.50%  '
